$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 310, pushing existing rows 310-337 down to 311-338
# (matches the canonical diff: dimension A1:R337 -> A1:R338, and every row's
# data from the old row N now lives in row N+1, with old row 337 becoming new row 338).
$ws.Rows(310).Insert()

# Populate the newly inserted row 310 with the new record's data.
$ws.Range("A310").Value = 10
$ws.Range("B310").Value = "Vega Modelo de Temuco"
$ws.Range("C310").Value = "La Araucanía"
$ws.Range("D310").Value = 44585
$ws.Range("E310").Value = 9
$ws.Range("F310").Value = 100114014
$ws.Range("G310").Value = "Betarraga"
$ws.Range("H310").Value = "Sin especificar"
$ws.Range("I310").Value = "Primera"
$ws.Range("J310").Value = 145
$ws.Range("K310").Value = 7000
$ws.Range("L310").Value = 8000
$ws.Range("M310").Value = 7448
$ws.Range("N310").Value = "$/docena de paquetes"
$ws.Range("O310").Value = "Provincia de Cautín"
$ws.Range("P310").Value = 621
$ws.Range("Q310").Value = 12
$ws.Range("R310").Value = "Hortaliza"
